# leet code 14; longest prefix
# Adds a new row (row 13) to the "July" worksheet for the
# "Longest common prefix" LeetCode problem.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for the new row ---
$ws.Range("A13").Value = "easy"
$ws.Range("B13").Value = 14
$ws.Range("C13").Value = "Longest common prefix"
$ws.Range("D13").Value = "string manipulation"
$ws.Range("E13").Value = "you an take any string in the list as a string to compare"
$ws.Range("F13").Value = "https://leetcode.com/problems/longest-common-prefix/submissions/1332877197 "

# --- Styling: match the look of the rest of the table ---
# Columns A, B, E use centered alignment (same as the rows above them).
$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("E13").HorizontalAlignment = -4108

# Row height matches the rest of the table.
$ws.Rows.Item(13).RowHeight = 17

# --- Hyperlink for the submission link in column F ---
$ws.Hyperlinks.Add($ws.Range("F13"), "https://leetcode.com/problems/longest-common-prefix/submissions/1332877197 ")

# Re-apply the hyperlink cell style (Hyperlinks.Add resets it), matching the
# look used by the other hyperlink cells (e.g. F5, which also uses the
# built-in "Hyperlink" cell style).
$ws.Range("F13").Style = "Hyperlink"

# --- Restore the selection that was active when the workbook was saved ---
[void]$ws.Range("E16").Select()
